$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full data set for rows 2-8 ("Mandataire" entries) after the "Rappel AV cloture" fix:
# - row order reshuffled (ZERNAKH ABDELLAH moved down, new rows inserted)
# - a couple of "N° de contrat" values / amounts updated
# - new beneficiaries added (MOHAMED BADRANE, NASIRI HASNAA) plus two extra
#   "NOUBAIL MOHAMMED" virements ("Rappel AV")
$data = @(
    @("NOUBAIL MOUNTASSIR","Q251990","007400000313200019604463","KHOURIBGA ZELLAKA","AWB","Direction régionale","905/TADLA OUARDIGHA ZAYANE","mensuelle",6750,675,6075),
    @("NOUBAIL MOHAMMED","IR801997","007400000313200019604463","KHOURIBGA ZELLAKA","AWB","Direction régionale","905/TADLA OUARDIGHA ZAYANE","mensuelle",6750,675,6075),
    @("MOHAMED BADRANE","I83603","225400000805987601012173","KHOURIBGA","CA","Point de vente","605/KHOURIBGA NAHDA","mensuelle",7500,375,7125),
    @("ZERNAKH ABDELLAH","IB19558","145101211406073828000084","MARRAKECH BENI MELLAL","BP","Point de vente","052/FKIH BEN SALEH/AV1","mensuelle",12000,0,12000),
    @("NASIRI HASNAA","","546576878798989898090090","","CIH","Logement de fonction","905/LF/TADLA OUARDIGHA ZAYANE","mensuelle",9999.99,999.99,9000),
    @("NOUBAIL MOHAMMED","IR801997","007400000313200019604463","KHOURIBGA ZELLAKA","AWB","Direction régionale","035/TES/AV1","mensuelle",1000,100,900),
    @("NOUBAIL MOHAMMED","IR801997","007400000313200019604463","KHOURIBGA ZELLAKA","AWB","Direction régionale","035/TES/AV1","mensuelle",4000,400,3600)
)

# Column C ("N° de compte") holds long numeric-looking account numbers with
# leading zeros -> must stay text, never become a rounded/scientific number.
$ws.Range("C2:C8").NumberFormat = "@"

$row = 2
foreach ($rowData in $data) {
    $col = 1
    foreach ($val in $rowData) {
        $ws.Cells.Item($row, $col).Value = $val
        $col++
    }
    $row++
}

# Totals row, now on row 9 (table grew from 3 to 7 detail rows)
$ws.Range("A9:H9").Value = " "
$ws.Cells.Item(9, 9).Value = 47999.99
$ws.Cells.Item(9, 10).Value = 3224.99
$ws.Cells.Item(9, 11).Value = 44775

$wb.Save()
